$p = $ppt.ActivePresentation

# Swap slide 2 ("Let the games begin!") and slide 3 ("What is Gamma?")
# by moving the slide currently at position 3 to position 2.
$s3 = $p.Slides.Item(3)
$s3.MoveTo(2)

Write-Host "done"
